$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the obsolete years 2004-2009 (old rows 2-7); remaining rows shift up,
# so former row 8 (2010) becomes row 2, ... former row 18 (2020) becomes row 12.
$ws.Rows("2:7").Delete()

# Step 2: refresh the precision of a handful of values for years 2010-2020 that were
# re-rounded (B/F/I to 4 decimals, D/G/J to whole numbers) in the source update.
$updates = @(
    @{ Row = 2;  B = 279.938;     D = 2802;  F = 487.5808;   G = 17021; I = 1268.0054;  J = 16503 },
    @{ Row = 3;  B = 84.725587;   D = 2676.33; F = 678.799723; G = 19027.8; I = 1165.832608; J = 16093.58 },
    @{ Row = 4;  B = 76.9686;     D = 2442.28; F = 795.0377;  G = 21207.53; I = 1114.8032; J = 15682.86 },
    @{ Row = 5;  B = 62.798877;   D = 1942.99; F = 900.990392; G = 23783.44; I = 1109.729791; J = 15101.95 },
    @{ Row = 6;  B = 55.951269;   D = 1757.01; F = 964.378286; G = 25972.94; I = 1082.849026; J = 14378.41 },
    @{ Row = 7;  B = 47.1378;     D = 1322;  F = 1040.7906;  G = 28561; I = 1039.2169;  J = 13955 },
    @{ Row = 8;  B = 44.0944;     D = 1085;  F = 1171.7186;  G = 30856; I = 1078.8042;  J = 13744 },
    @{ Row = 9;  B = 27.0882;     D = 752;   F = 1263.7546;  G = 33934; I = 998.8088;   J = 12616 },
    @{ Row = 10; B = 29.7893;     D = 779;   F = 1443.9538;  G = 36902; I = 1015.3298;  J = 11782 },
    @{ Row = 11; B = 27.6841;     D = 675;   F = 1608.557;   G = 39025; I = 1040.811;   J = 11297 },
    @{ Row = 12; B = 23.1447;     D = 548;   F = 1563.702;   G = 41302; I = 833.7109;   J = 10767 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 9).Value = $u.I
    $ws.Cells.Item($r, 10).Value = $u.J
}

# Step 3: append the new 2021 row (row 13), matching the formatting of column A's
# existing year cells (bold, centered/top aligned, thin border).
$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value = "2021年"
$ws.Cells.Item(13, 2).Value = 18.723418
$ws.Cells.Item(13, 3).Value = 969.97
$ws.Cells.Item(13, 4).Value = 455.77
$ws.Cells.Item(13, 5).Value = 9165.01
$ws.Cells.Item(13, 6).Value = 1721.061247
$ws.Cells.Item(13, 7).Value = 44195.53
$ws.Cells.Item(13, 8).Value = 929087.71
$ws.Cells.Item(13, 9).Value = 860.684058
$ws.Cells.Item(13, 10).Value = 10180.49
$ws.Cells.Item(13, 11).Value = 2909.98

Write-Output "done"
